# Updates cryptocurrency price/symbol data in the "cryptos" worksheet to
# reflect the latest GitHub Actions scrape (commit: "Updated symbol list on
# Sun Dec 18 17:44:34 UTC 2022 with GitHub Actions").
#
# Column D holds numeric-looking price text that must stay stored as TEXT
# (as it was originally), so Set-TextValue forces the cell to Text format
# before assigning the value -- this prevents Excel from silently
# re-typing it as a number.

function Set-TextValue {
    param($ws, $cellRef, $val)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "247.19"
Set-TextValue $ws "D3" "22.42"
Set-TextValue $ws "D4" "5.475"
Set-TextValue $ws "D5" "0.05618"
Set-TextValue $ws "D6" "6.458"
Set-TextValue $ws "D7" "0.8042"
Set-TextValue $ws "D8" "1.039"
Set-TextValue $ws "D9" "0.1422"
Set-TextValue $ws "D10" "0.07318"
Set-TextValue $ws "D11" "0.03171"
Set-TextValue $ws "D12" "0.02920"
Set-TextValue $ws "D14" "0.001672"
Set-TextValue $ws "D15" "3.222"
Set-TextValue $ws "D16" "0.04736"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws "D17" "0.006468"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws "D18" "0.005068"
$ws.Range("E18").Value = "17HotbitTokenHTB"
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws "D19" "0.001051"
$ws.Range("E19").Value = "18BitKanKAN"
$ws.Range("B20").Value = "NitroEx"
$ws.Range("C20").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws "D20" "0.0001503"
$ws.Range("E20").Value = "19NitroExNTX"
$ws.Range("B21").Value = "LEO"
$ws.Range("C21").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws "D21" "3.985"
$ws.Range("E21").Value = "20LEOLEO"
$ws.Range("B22").Value = "GateToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws "D22" "3.380"
$ws.Range("E22").Value = "21GateTokenGT"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws "D23" "2.102"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D24" "0.01166"
$ws.Range("E24").Value = "23OneONEBestin24h"
$ws.Range("E26").Value = "25ProBitTokenPROB"
Set-TextValue $ws "D27" "0.0003306"
Set-TextValue $ws "D40" "0.04174"
Set-TextValue $ws "D41" "0.006873"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003507"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1038"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue $ws "D44" "0.009099"
Set-TextValue $ws "D45" "0.00005661"
Set-TextValue $ws "D47" "0.6813"
Set-TextValue $ws "D48" "0.01555"

Write-Host "Applied symbol list update."
